# Apply trade #32 close: update Summary, Strategy Status, and append the
# new trade row to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.76      # Current Capital
$summary.Range("B4").Value = -0.24        # Total P&L $
$summary.Range("B6").Value = 32           # Total Trades
$summary.Range("B8").Value = 17           # Losing Trades
$summary.Range("B9").Value = 28.12        # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.76000000000001   # Capital
$status.Range("D4").Value = 32                  # Trades
$status.Range("E4").Value = -0.24               # P&L $
$status.Range("F4").Value = -0.24               # P&L %
$status.Range("G4").Value = 28.12               # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the new trade (#32) to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------
# Note: the Date column ("2026-02-17") must stay plain text, so it is
# entered with a leading apostrophe to stop Excel auto-converting it to
# a date serial number (same effect as typing it directly into a cell).
$tradeRow = @(32, "'2026-02-17", "15:22:39", "MarketMaking", "DOWN", 0.26, 0.25017, "CLOSED", -3.7807, -0.01, 99.76000000000001, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = 33
    for ($i = 0; $i -lt $tradeRow.Length; $i++) {
        $ws.Cells.Item($newRow, $i + 1).Value = $tradeRow[$i]
    }
}

$wb.Save()
